$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: split Timestamp into TimestampStart / TimestampEnd,
# and move AnswerIsCorrect to new column I
$ws.Range("G1").Value = "TimestampStart"
$ws.Range("H1").Value = "TimestampEnd"
$ws.Range("I1").Value = "AnswerIsCorrect"

# Remove the custom width on column H (was set to 16 for "Timestamp")
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Clear out the old empty placeholder rows (11-20) that were left over
$ws.Rows("11:20").Delete()

# Set the active selection cell like in the new sheet
$ws.Range("E14").Select()
